# issue #5: add legislator_id, name, date into dataframe
#
# The stock holdings sheet ("股票") gains three new trailing columns:
#   date (H), legislator_name (I), legislator_id (J)
# populated for every existing data row with the report's filing date,
# the legislator's name and numeric id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Last populated data row (header is row 1, data starts row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# --- headers (row 1) -------------------------------------------------
# Re-use the existing header formatting (bold / bordered) from column G.
$ws.Cells.Item(1, 7).Copy() | Out-Null
$ws.Range($ws.Cells.Item(1, 8), $ws.Cells.Item(1, 10)).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# --- data rows ---------------------------------------------------------
# Re-use the plain data formatting from column G for the new data cells.
$ws.Cells.Item(2, 7).Copy() | Out-Null
$ws.Range($ws.Cells.Item(2, 8), $ws.Cells.Item($lastRow, 10)).PasteSpecial(-4122) | Out-Null

# The date column must stay plain text ("2011-12-29"), not be reinterpreted
# as a date serial number, so force a text number format on it first.
$ws.Range($ws.Cells.Item(2, 8), $ws.Cells.Item($lastRow, 8)).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "2011-12-29"
    $ws.Cells.Item($r, 9).Value = "蔡正元"
    $ws.Cells.Item($r, 10).Value = 966
}

$ws.Application.CutCopyMode = $false
